# Generate Report for handoff
# The file "ef91d88f-832a-4db3-a242-cc78104128ea.md" is now ready for handoff again
# (status changes from "Handed back: in sync with en-US" to "Ready for handoff"),
# and the Latest Handoff Datetime is refreshed for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the ef91d88f-832a-4db3-a242-cc78104128ea.md file.
# Both the zh-cn (B) and de-de (C) status columns flip to "Ready for handoff".
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 status -> "Ready for handoff",
# and the Latest Handoff Datetime (column D) is updated.
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-01-18 03:41:30"

# de-de sheet: row 3 status -> "Ready for handoff",
# and the Latest Handoff Datetime (column D) is updated.
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-01-18 03:41:47"
